$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply hyperlink style (same as other B column cells) before setting the value
$ws.Range("B13").Style = "Hipervínculo"

# New row 13 data: URL link text in B13, description in C13
$ws.Range("B13").Value = "https://www.plasmic.app/blog/mastering-css-flexbox-with-plasmic"
$ws.Range("C13").Value = "Ejemplos de flexbox (imágenes y código)"

# Add the hyperlink itself
$ws.Hyperlinks.Add($ws.Range("B13"), "https://www.plasmic.app/blog/mastering-css-flexbox-with-plasmic")

# Re-assert the hyperlink style so B13 keeps sharing the same style index as B3:B12
$ws.Range("B13").Style = "Hipervínculo"

# Update selection to C14
$ws.Range("C14").Select()
